$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Test Id and Test case name for the Cart Module defect (column B)
$ws.Range("B3").Value = "TC_011"
$ws.Range("B4").Value = "TC_Cart_011"

# Reflect the new selection left by the edit
$ws.Range("B4").Select()
